# "fase de atacar, botón de pasar turno, mostrar informacion de tableros y cartas"
#
# Sheet "Hoja1":
#  - Remove the "Devolver estado mana" row from the "Fase 1 (por jugador)" block.
#  - Remove the "Devolver estado mana" row from the "Fases juego (por cada
#    jugador)" repeat block.
#  - Rename "Atacar" -> "Atacar ** pendiente comprobar maná suficiente" and
#    highlight it with a new fill colour (pending/in-progress marker).
#  - Add a new "Defender" step right below the (renamed) "Atacar" step.
#  - Update the selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Drop the first "Devolver estado mana" row (Fase 1 / por jugador block).
$ws.Rows.Item(17).Delete()

# Drop the second "Devolver estado mana" row (Fases juego / repeat block).
$ws.Rows.Item(22).Delete()

# Insert a new row for "Defender" right after the "Atacar" row, matching the
# plain "marker" style already used elsewhere in the sheet (e.g. C9).
$ws.Rows.Item(24).Insert()
$ws.Range("C24").Value = "Defender"
$ws.Range("C24").Style = $ws.Range("C9").Style

# "Atacar" is not fully implemented yet -> note it in the label and mark the
# cell with a new (orange) fill colour to flag it as pending work.
$ws.Range("C23").Value = "Atacar ** pendiente comprobar maná suficiente"
$ws.Range("C23").Interior.ThemeColor = 6

# Reflect the last-used selection when the sheet was saved.
$ws.Range("G22").Select()
